# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Terminal Hortofrutícola Agro Chillán - Mango"
# above the current row 29, shifting all subsequent rows down by 2
# (dimension grows from A1:T60 to A1:T62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 29-30; existing rows 29-60 become 31-62.
$ws.Rows("29:30").Insert()

# New row 29
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44482
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100108
$ws.Range("H29").Value = "Tropicales y subtropicales"
$ws.Range("I29").Value = 100108002
$ws.Range("J29").Value = "Mango"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 60
$ws.Range("N29").Value = 8500
$ws.Range("O29").Value = 9000
$ws.Range("P29").Value = 8750
$ws.Range("Q29").Value = "$/bandeja 4 kilos"
$ws.Range("R29").Value = "Brasil"
$ws.Range("S29").Value = 2188
$ws.Range("T29").Value = 4

# New row 30
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44482
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108002
$ws.Range("J30").Value = "Mango"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 7500
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 7750
$ws.Range("Q30").Value = "$/bandeja 4 kilos"
$ws.Range("R30").Value = "Perú"
$ws.Range("S30").Value = 1938
$ws.Range("T30").Value = 4
